$d = $word.ActiveDocument

# --- Create the three new character styles -------------------------------

$ganStyle = $d.Styles.Add("GaNStyle", 2)
$ganStyle.Font.Name = "Calibri"
$ganStyle.Font.Size = 14

$ganParagraph = $d.Styles.Add("GaNParagraph", 2)
$ganParagraph.Font.Name = "Calibri"
$ganParagraph.Font.Size = 10

$ganLinks = $d.Styles.Add("GaNLinks", 2)
$ganLinks.Font.Name = "Calibri"
$ganLinks.Font.Bold = $true
$ganLinks.Font.Color = 8388608
$ganLinks.Font.Size = 9.5
$ganLinks.Font.Underline = 1

# --- Apply GaNStyle to every "Dates à utiliser..." run -------------------

$rng = $d.Content
while ($rng.Find.Execute("Dates à utiliser pour la Campagne 2022 Constellation de Bootes: 14-23 mai, 13-22 juin, 12-21 juillet", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng.Style = "GaNStyle"
    $rng.Collapse(0)
}

# --- Apply GaNParagraph to the "Vous allez participer..." run ------------

$rng = $d.Content
while ($rng.Find.Execute("Vous allez participer à une campagne mondiale d’observation pour détecter les plus faibles étoiles visibles afin de mesurer la pollution lumineuse sur un site donné. Partout dans le monde, en localisant et en observant la Constellation de Bootes dans le ciel nocturne et en la comparant aux cartes stellaires, les participants, apprendront comment l’éclairage, dans leur environnement local, influence la pollution lumineuse. Vos contributions à la base de données en ligne permettront de mesurer la qualité du ciel nocturne.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng.Style = "GaNParagraph"
    $rng.Collapse(0)
}

# --- Apply GaNLinks to the "Les cartes figurant..." run ------------------

$rng = $d.Content
while ($rng.Find.Execute("Les cartes figurant dans ce document ont été établies par Jenik Hollan, CzechGlobe ((http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/).", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng.Style = "GaNLinks"
    $rng.Collapse(0)
}
